$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: a cell known to carry the plain default style (no explicit style index),
# used as a template to strip any stray style picked up while forcing text values.
$plainStyleTemplate = $ws.Range("B2")

# D2: "29.423.42" -> "29.426.18"
$ws.Range("D2").Value = "29.426.18"

# E2: "  +0.61%  " -> "  +0.56%  "
$ws.Range("E2").Value = "  +0.56%  "

# D3: "1.875.88" -> "1.876.26"
$ws.Range("D3").Value = "1.876.26"

# E3: "  +1.05%  " -> "  +1.03%  "
$ws.Range("E3").Value = "  +1.03%  "

# D4: "1.001" -> "1.000" (force text to avoid numeric auto-conversion)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = $plainStyleTemplate.Style

# D5: "0.7122" -> "0.7127" (force text to avoid numeric auto-conversion)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7127"
$ws.Range("D5").Style = $plainStyleTemplate.Style

# E5: "  +1.50%  " -> "  +1.71%  "
$ws.Range("E5").Value = "  +1.71%  "

# D6: "241.79" -> "241.65" (force text to avoid numeric auto-conversion)
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.65"
$ws.Range("D6").Style = $plainStyleTemplate.Style

# E6: "  +1.63%  " -> "  +1.51%  "
$ws.Range("E6").Value = "  +1.51%  "

# D8: "0.07847" -> "0.07831" (force text to avoid numeric auto-conversion)
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07831"
$ws.Range("D8").Style = $plainStyleTemplate.Style

# E8: "  -2.35%  " -> "  -2.47%  "
$ws.Range("E8").Value = "  -2.47%  "

# D9: "0.3108" -> "0.3106" (force text to avoid numeric auto-conversion)
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3106"
$ws.Range("D9").Style = $plainStyleTemplate.Style

# E9: "  +2.97%  " -> "  +2.83%  "
$ws.Range("E9").Value = "  +2.83%  "

# D10: "25.17" -> "25.14" (force text to avoid numeric auto-conversion)
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "25.14"
$ws.Range("D10").Style = $plainStyleTemplate.Style

# E10: "  +7.18%  " -> "  +7.01%  "
$ws.Range("E10").Value = "  +7.01%  "

# D11: "0.08242" -> "0.08235" (force text to avoid numeric auto-conversion)
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08235"
$ws.Range("D11").Style = $plainStyleTemplate.Style

# E11: "  +0.73%  " -> "  +0.57%  "
$ws.Range("E11").Value = "  +0.57%  "

# D12: "0.7273" -> "0.7271" (force text to avoid numeric auto-conversion)
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7271"
$ws.Range("D12").Style = $plainStyleTemplate.Style

# E12: "  +2.91%  " -> "  +2.85%  "
$ws.Range("E12").Value = "  +2.85%  "

# B13: "Polkadot" -> "WrappedEther"
$ws.Range("B13").Value = "WrappedEther"

# C13: "https://coinranking.com/coin/25W7FG7om+polkadot-dot" -> "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"

# D13: "5.258" -> "1.868.80"
$ws.Range("D13").Value = "1.868.80"

# E13: "  +1.06%  " -> "  +0.18%  "
$ws.Range("E13").Value = "  +0.18%  "

# B14: "WrappedEther" -> "Polkadot"
$ws.Range("B14").Value = "Polkadot"

# C14: "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth" -> "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"

# D14: "1.855.80" -> "5.256" (force text to avoid numeric auto-conversion)
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.256"
$ws.Range("D14").Style = $plainStyleTemplate.Style

# E14: "  -0.69%  " -> "  +0.97%  "
$ws.Range("E14").Value = "  +0.97%  "

# D15: "90.77" -> "90.78" (force text to avoid numeric auto-conversion)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.78"
$ws.Range("D15").Style = $plainStyleTemplate.Style

# E15: "  +1.27%  " -> "  +1.22%  "
$ws.Range("E15").Value = "  +1.22%  "

# D16: "29.422.73" -> "29.428.84"
$ws.Range("D16").Value = "29.428.84"

# E16: "  +0.39%  " -> "  +0.45%  "
$ws.Range("E16").Value = "  +0.45%  "

# D17: "5.902" -> "5.901" (force text to avoid numeric auto-conversion)
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.901"
$ws.Range("D17").Style = $plainStyleTemplate.Style

# E17: "  +1.17%  " -> "  +1.26%  "
$ws.Range("E17").Value = "  +1.26%  "

# D18: "247.25" -> "247.32" (force text to avoid numeric auto-conversion)
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "247.32"
$ws.Range("D18").Style = $plainStyleTemplate.Style

# E18: "  +4.05%  " -> "  +4.08%  "
$ws.Range("E18").Value = "  +4.08%  "

# D19: "0.000007861" -> "0.000007864" (force text to avoid numeric auto-conversion)
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007864"
$ws.Range("D19").Style = $plainStyleTemplate.Style

# E19: "  -0.44%  " -> "  -0.37%  "
$ws.Range("E19").Value = "  -0.37%  "

# E20: "  -0.13%  " -> "  -0.24%  "
$ws.Range("E20").Value = "  -0.24%  "

# D21: "2.115.17" -> "2.117.40"
$ws.Range("D21").Value = "2.117.40"

# E21: "  -0.32%  " -> "  -0.06%  "
$ws.Range("E21").Value = "  -0.06%  "

# D22: "0.9998" -> "0.9992" (force text to avoid numeric auto-conversion)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9992"
$ws.Range("D22").Style = $plainStyleTemplate.Style

# E22: "  +0.04%  " -> "  -0.18%  "
$ws.Range("E22").Value = "  -0.18%  "

# D23: "7.962" -> "7.977" (force text to avoid numeric auto-conversion)
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.977"
$ws.Range("D23").Style = $plainStyleTemplate.Style

# E23: "  +6.93%  " -> "  +7.15%  "
$ws.Range("E23").Value = "  +7.15%  "

# D24: "1.000" -> "0.9995" (force text to avoid numeric auto-conversion)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9995"
$ws.Range("D24").Style = $plainStyleTemplate.Style

# E24: "  -0.06%  " -> "  -0.20%  "
$ws.Range("E24").Value = "  -0.20%  "

# D25: "0.1575" -> "0.1574" (force text to avoid numeric auto-conversion)
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1574"
$ws.Range("D25").Style = $plainStyleTemplate.Style

# E25: "  +10.33%  " -> "  +10.46%  "
$ws.Range("E25").Value = "  +10.46%  "

# D26: "163.55" -> "163.79" (force text to avoid numeric auto-conversion)
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.79"
$ws.Range("D26").Style = $plainStyleTemplate.Style

# E26: "  +0.39%  " -> "  +0.41%  "
$ws.Range("E26").Value = "  +0.41%  "

# D27: "8.983" -> "8.995" (force text to avoid numeric auto-conversion)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.995"
$ws.Range("D27").Style = $plainStyleTemplate.Style

# D28: "18.28" -> "18.27" (force text to avoid numeric auto-conversion)
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.27"
$ws.Range("D28").Style = $plainStyleTemplate.Style

# E28: "  +1.06%  " -> "  +0.96%  "
$ws.Range("E28").Value = "  +0.96%  "

# E29: "  -3.42%  " -> "  -3.43%  "
$ws.Range("E29").Value = "  -3.43%  "

# D30: "1.490" -> "1.495" (force text to avoid numeric auto-conversion)
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.495"
$ws.Range("D30").Style = $plainStyleTemplate.Style

# E30: "  +0.99%  " -> "  +1.48%  "
$ws.Range("E30").Value = "  +1.48%  "

# D31: "4.359" -> "4.360" (force text to avoid numeric auto-conversion)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.360"
$ws.Range("D31").Style = $plainStyleTemplate.Style

# E31: "  -0.39%  " -> "  -0.35%  "
$ws.Range("E31").Value = "  -0.35%  "

# D32: "4.122" -> "4.124" (force text to avoid numeric auto-conversion)
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.124"
$ws.Range("D32").Style = $plainStyleTemplate.Style

# E32: "  +2.24%  " -> "  +2.13%  "
$ws.Range("E32").Value = "  +2.13%  "

# D33: "0.05300" -> "0.05306" (force text to avoid numeric auto-conversion)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05306"
$ws.Range("D33").Style = $plainStyleTemplate.Style

# E34: "  +0.36%  " -> "  +0.40%  "
$ws.Range("E34").Value = "  +0.40%  "

# E35: "  +3.11%  " -> "  +2.98%  "
$ws.Range("E35").Value = "  +2.98%  "

# D36: "0.7212" -> "0.7220" (force text to avoid numeric auto-conversion)
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7220"
$ws.Range("D36").Style = $plainStyleTemplate.Style

# E36: "  +0.32%  " -> "  +0.13%  "
$ws.Range("E36").Value = "  +0.13%  "

# E37: "  -0.55%  " -> "  -0.90%  "
$ws.Range("E37").Value = "  -0.90%  "

# D38: "0.01860" -> "0.01859" (force text to avoid numeric auto-conversion)
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01859"
$ws.Range("D38").Style = $plainStyleTemplate.Style

# E38: "  +0.59%  " -> "  +0.44%  "
$ws.Range("E38").Value = "  +0.44%  "

# D39: "1.249.58" -> "1.244.76"
$ws.Range("D39").Value = "1.244.76"

# E39: "  +8.26%  " -> "  +7.95%  "
$ws.Range("E39").Value = "  +7.95%  "

# D40: "2.722" -> "2.727" (force text to avoid numeric auto-conversion)
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.727"
$ws.Range("D40").Style = $plainStyleTemplate.Style

# E40: "  -0.04%  " -> "  +0.13%  "
$ws.Range("E40").Value = "  +0.13%  "

# D41: "0.9081" -> "0.9077" (force text to avoid numeric auto-conversion)
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9077"
$ws.Range("D41").Style = $plainStyleTemplate.Style

# E41: "  -3.22%  " -> "  -3.33%  "
$ws.Range("E41").Value = "  -3.33%  "

# D42: "73.68" -> "73.82" (force text to avoid numeric auto-conversion)
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.82"
$ws.Range("D42").Style = $plainStyleTemplate.Style

# E42: "  +4.57%  " -> "  +4.78%  "
$ws.Range("E42").Value = "  +4.78%  "

# D43: "6.141" -> "6.146" (force text to avoid numeric auto-conversion)
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.146"
$ws.Range("D43").Style = $plainStyleTemplate.Style

# D45: "103.19" -> "103.18" (force text to avoid numeric auto-conversion)
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "103.18"
$ws.Range("D45").Style = $plainStyleTemplate.Style

# E45: "  +0.21%  " -> "  +0.30%  "
$ws.Range("E45").Value = "  +0.30%  "

# D46: "0.5330" -> "0.5324" (force text to avoid numeric auto-conversion)
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5324"
$ws.Range("D46").Style = $plainStyleTemplate.Style

# E46: "  +0.69%  " -> "  +0.66%  "
$ws.Range("E46").Value = "  +0.66%  "

# D47: "1.766" -> "1.770" (force text to avoid numeric auto-conversion)
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.770"
$ws.Range("D47").Style = $plainStyleTemplate.Style

# E47: "  +0.61%  " -> "  +0.68%  "
$ws.Range("E47").Value = "  +0.68%  "

# D48: "2.921" -> "2.930" (force text to avoid numeric auto-conversion)
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.930"
$ws.Range("D48").Style = $plainStyleTemplate.Style

# E48: "  +13.08%  " -> "  +12.96%  "
$ws.Range("E48").Value = "  +12.96%  "

# E49: "  -0.03%  " -> "  -0.32%  "
$ws.Range("E49").Value = "  -0.32%  "

# D50: "0.4311" -> "0.4313" (force text to avoid numeric auto-conversion)
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4313"
$ws.Range("D50").Style = $plainStyleTemplate.Style

# E50: "  +1.20%  " -> "  +1.14%  "
$ws.Range("E50").Value = "  +1.14%  "

# E51: "  +0.85%  " -> "  +0.84%  "
$ws.Range("E51").Value = "  +0.84%  "
